$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "true"/"false" status values (column F) with "ativo"/"inativo"
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $c = $ws.Cells.Item($r, 6)
    if ($c.Value2 -eq $true) {
        $c.Value = "ativo"
    } elseif ($c.Value2 -eq $false) {
        $c.Value = "inativo"
    }
}

# Add a new empty, underlined cell at M16
$m16 = $ws.Range("M16")
$m16.Font.Underline = $true

# Update the selection to M16, matching the saved workbook state
$m16.Select()

# Configure page setup (paper size + orientation) for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
